# Re-pulled data: update dSF ("F" column) values for the affected rows.
# Mirrors the new closing-spread data while leaving dS0 ("E" column) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = -6
    4  = 9
    7  = 9
    11 = -1
    12 = 0
    14 = -2
    19 = -2
    20 = -1
    21 = 3
    27 = -4
    32 = -4
    35 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
